# Actualiza la base de datos del Estado de Cuenta:
# elimina los periodos de mora anteriores y agrega los nuevos periodos
# (2101-2107) para cada trabajador, reordenando la tabla de detalle
# (filas 16 a 36) en el orden: JAIR (2008), JESUS DAVID (2107..2101,2008),
# OLISMAIDA (2107..2101), VANESSA (2107..2103).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

$rows = @(
    @{Row=16; C="1047457290"; D="JAIR JOSE RACERO BANQUEZ"; E="2008"; F=35112},
    @{Row=17; C="1143367910"; D="JESUS DAVID LOPEZ BROME"; E="2107"; F=29260},
    @{Row=18; C="1143367910"; D="JESUS DAVID LOPEZ BROME"; E="2106"; F=35112},
    @{Row=19; C="1143367910"; D="JESUS DAVID LOPEZ BROME"; E="2105"; F=35112},
    @{Row=20; C="1143367910"; D="JESUS DAVID LOPEZ BROME"; E="2104"; F=35112},
    @{Row=21; C="1143367910"; D="JESUS DAVID LOPEZ BROME"; E="2103"; F=35112},
    @{Row=22; C="1143367910"; D="JESUS DAVID LOPEZ BROME"; E="2102"; F=35112},
    @{Row=23; C="1143367910"; D="JESUS DAVID LOPEZ BROME"; E="2101"; F=35112},
    @{Row=24; C="1143367910"; D="JESUS DAVID LOPEZ BROME"; E="2008"; F=35112},
    @{Row=25; C="1201219362"; D="OLISMAIDA GARCIA SOTO"; E="2107"; F=29260},
    @{Row=26; C="1201219362"; D="OLISMAIDA GARCIA SOTO"; E="2106"; F=35112},
    @{Row=27; C="1201219362"; D="OLISMAIDA GARCIA SOTO"; E="2105"; F=35112},
    @{Row=28; C="1201219362"; D="OLISMAIDA GARCIA SOTO"; E="2104"; F=35112},
    @{Row=29; C="1201219362"; D="OLISMAIDA GARCIA SOTO"; E="2103"; F=35112},
    @{Row=30; C="1201219362"; D="OLISMAIDA GARCIA SOTO"; E="2102"; F=35112},
    @{Row=31; C="1201219362"; D="OLISMAIDA GARCIA SOTO"; E="2101"; F=35112},
    @{Row=32; C="1047424362"; D="VANESSA OSORIO SIMANCAS"; E="2107"; F=29260},
    @{Row=33; C="1047424362"; D="VANESSA OSORIO SIMANCAS"; E="2106"; F=35112},
    @{Row=34; C="1047424362"; D="VANESSA OSORIO SIMANCAS"; E="2105"; F=35112},
    @{Row=35; C="1047424362"; D="VANESSA OSORIO SIMANCAS"; E="2104"; F=35112},
    @{Row=36; C="1047424362"; D="VANESSA OSORIO SIMANCAS"; E="2103"; F=35112}
)

foreach ($r in $rows) {
    $ws.Cells.Item($r.Row, 3).Value = $r.C
    $ws.Cells.Item($r.Row, 4).Value = $r.D
    $ws.Cells.Item($r.Row, 5).Value = $r.E
    $ws.Cells.Item($r.Row, 6).Value = $r.F
}
